# ---------------------------------------------------------------------------
# Adds a "2022-Q3" quarter to the workbook:
#   1. "总计" (summary) sheet gets a new row 2 with the 2022-Q3 totals; all
#      the other quarter rows shift down by one and their running index
#      (column A) is bumped by one.
#   2. A brand-new worksheet named "2022-Q3", placed right after "总计",
#      holding the per-fund holding breakdown for that quarter.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper sheet references
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Make sure the summary sheet is the active one before we start adding
# sheets - a freshly added sheet otherwise shadows whatever sheet happens
# to be active already.
$summary.Activate()

# ---------------------------------------------------------------------------
# 1) Update "总计" - insert the 2022-Q3 row at the top of the data block.
# ---------------------------------------------------------------------------
$summary.Rows("2:2").Insert()

# The freshly inserted row picks up stray formatting from the Insert() -
# strip it so text/number cells go back to "no explicit style" like the
# rest of the sheet.
$summary.Range("B2:D2").ClearFormats()

$summaryRows = @(
    @{ Row = 2; Idx = 0; Period = "2022-Q3"; Count = 21; Value = 2.32 },
    @{ Row = 3; Idx = 1; Period = "2022-Q2"; Count = 11; Value = 1.59 },
    @{ Row = 4; Idx = 2; Period = "2022-Q1"; Count = 22; Value = 2.41 },
    @{ Row = 5; Idx = 3; Period = "2021-Q4"; Count = 7;  Value = 1.3 },
    @{ Row = 6; Idx = 4; Period = "2021-Q3"; Count = 4;  Value = 0.7 },
    @{ Row = 7; Idx = 5; Period = "2021-Q2"; Count = 5;  Value = 0.64 },
    @{ Row = 8; Idx = 6; Period = "2021-Q1"; Count = 10; Value = 0.67 },
    @{ Row = 9; Idx = 7; Period = "2020-Q4"; Count = 5;  Value = 0.9 }
)

foreach ($r in $summaryRows) {
    $rowNum = $r.Row
    $summary.Range("A$rowNum").Value = $r.Idx
    $summary.Range("B$rowNum").Value = $r.Period
    $summary.Range("C$rowNum").Value = $r.Count
    $summary.Range("D$rowNum").Value = $r.Value
}

# Re-apply the "index column" formatting (bold / bordered / centered) that
# the rest of column A already carries. Row 3 (the shifted former row 2)
# still has the original style, so copy from there onto the whole column.
$summary.Range("A3").Copy() | Out-Null
$summary.Range("A2:A9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Create the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = [char]([int][char]'B' + $i)
    $q3.Range("$col`1").Value = $headers[$i]
}

$funds = @(
    @{ Idx = 0;  Code = "002621"; Name = "中欧消费主题股票A";                 Size = "10.65"; Stock = "88.44"; Pct = "5.54"; Value = "0.5900"; Rank = 10 },
    @{ Idx = 1;  Code = "003751"; Name = "万家瑞隆混合A";                     Size = "14.54"; Stock = "93.96"; Pct = "2.65"; Value = "0.3853"; Rank = 10 },
    @{ Idx = 2;  Code = "002697"; Name = "中欧消费主题股票C";                 Size = "5.37";  Stock = "88.44"; Pct = "5.54"; Value = "0.2975"; Rank = 10 },
    @{ Idx = 3;  Code = "001195"; Name = "工银农业产业股票";                  Size = "5.53";  Stock = "80.67"; Pct = "4.10"; Value = "0.2267"; Rank = 4 },
    @{ Idx = 4;  Code = "010852"; Name = "中欧内需成长混合A";                 Size = "3.70";  Stock = "90.11"; Pct = "5.66"; Value = "0.2094"; Rank = 7 },
    @{ Idx = 5;  Code = "970023"; Name = "天风天盈一年定期开放混合";          Size = "1.24";  Stock = "70.03"; Pct = "9.42"; Value = "0.1168"; Rank = 1 },
    @{ Idx = 6;  Code = "012442"; Name = "永赢稳健增长一年持有期混合E";       Size = "9.56";  Stock = "26.04"; Pct = "0.98"; Value = "0.0937"; Rank = 10 },
    @{ Idx = 7;  Code = "005620"; Name = "中欧品质消费股票A";                 Size = "1.62";  Stock = "91.58"; Pct = "5.77"; Value = "0.0935"; Rank = 7 },
    @{ Idx = 8;  Code = "005526"; Name = "工银瑞信新生代消费灵活配置混合";    Size = "1.95";  Stock = "94.51"; Pct = "2.90"; Value = "0.0566"; Rank = 10 },
    @{ Idx = 9;  Code = "005621"; Name = "中欧品质消费股票C";                 Size = "0.97";  Stock = "91.58"; Pct = "5.77"; Value = "0.0560"; Rank = 7 },
    @{ Idx = 10; Code = "009932"; Name = "永赢稳健增长一年持有期混合A";       Size = "4.93";  Stock = "26.04"; Pct = "0.98"; Value = "0.0483"; Rank = 10 },
    @{ Idx = 11; Code = "015384"; Name = "万家瑞隆混合C";                     Size = "1.52";  Stock = "93.96"; Pct = "2.65"; Value = "0.0403"; Rank = 10 },
    @{ Idx = 12; Code = "010853"; Name = "中欧内需成长混合C";                 Size = "0.54";  Stock = "90.11"; Pct = "5.66"; Value = "0.0306"; Rank = 7 },
    @{ Idx = 13; Code = "011536"; Name = "惠升惠益混合A";                     Size = "0.64";  Stock = "21.48"; Pct = "2.55"; Value = "0.0163"; Rank = 5 },
    @{ Idx = 14; Code = "008491"; Name = "万家周期优势企业混合A";             Size = "0.61";  Stock = "93.50"; Pct = "2.66"; Value = "0.0162"; Rank = 8 },
    @{ Idx = 15; Code = "005599"; Name = "汇安量化优选灵活配置混合A";         Size = "0.51";  Stock = "94.40"; Pct = "2.73"; Value = "0.0139"; Rank = 10 },
    @{ Idx = 16; Code = "009128"; Name = "明亚价值长青混合A";                 Size = "0.38";  Stock = "57.73"; Pct = "3.51"; Value = "0.0133"; Rank = 5 },
    @{ Idx = 17; Code = "011537"; Name = "惠升惠益混合C";                     Size = "0.32";  Stock = "21.48"; Pct = "2.55"; Value = "0.0082"; Rank = 5 },
    @{ Idx = 18; Code = "008492"; Name = "万家周期优势企业混合C";             Size = "0.14";  Stock = "93.50"; Pct = "2.66"; Value = "0.0037"; Rank = 8 },
    @{ Idx = 19; Code = "005600"; Name = "汇安量化优选灵活配置混合C";         Size = "0.02";  Stock = "94.40"; Pct = "2.73"; Value = "0.0005"; Rank = 10 },
    @{ Idx = 20; Code = "009129"; Name = "明亚价值长青混合C";                 Size = "0.00";  Stock = "57.73"; Pct = "3.51"; Value = $null;     Rank = 5 }
)

# Text-typed columns (matching the source file, which stores these
# numeric-looking values as text): B, C, D, E, F, G (G22 is the one
# genuine numeric exception, handled separately below).
$lastRow = 1 + $funds.Length
$q3.Range("B2:G$lastRow").NumberFormat = "@"

$row = 2
foreach ($f in $funds) {
    $q3.Range("A$row").Value = $f.Idx
    $q3.Range("B$row").Value = $f.Code
    $q3.Range("C$row").Value = $f.Name
    $q3.Range("D$row").Value = $f.Size
    $q3.Range("E$row").Value = $f.Stock
    $q3.Range("F$row").Value = $f.Pct

    if ($null -ne $f.Value) {
        $q3.Range("G$row").Value = $f.Value
    }

    $q3.Range("H$row").Value = $f.Rank

    $row++
}

# Drop the temporary text formatting everywhere except G22, which Excel
# stores as a genuine number (0) in the source data.
$q3.Range("B2:G$lastRow").ClearFormats()
$q3.Range("G$lastRow").Value = 0

# ---------------------------------------------------------------------------
# Match formatting: header row + index column use the bold/bordered/centered
# style ("s=2" in the original file) that's already present on "总计".
# ---------------------------------------------------------------------------
$summary.Range("B1").Copy() | Out-Null
$q3.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$summary.Range("A3").Copy() | Out-Null
$q3.Range("A2:A$lastRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$q3.Range("A1").Select()

# Restore the original active-sheet state ("总计" was the active tab).
$summary.Activate()
$summary.Range("A1").Select()
